$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '75.926.68'
$ws.Range('E2').Value = '  +9.19%  '
$ws.Range('D3').Value = '2.690.87'
$ws.Range('E3').Value = '  +11.00%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '187.75'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +12.88%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '589.53'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.45%  '
$ws.Range('E8').Value = '  +5.22%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.195'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +15.67%  '
$ws.Range('D10').Value = '2.690.14'
$ws.Range('E10').Value = '  +10.98%  '
$ws.Range('E11').Value = '  +1.52%  '
$ws.Range('E13').Value = '  +2.01%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.189.21'
$ws.Range('E14').Value = '  +11.07%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '75.720.52'
$ws.Range('E15').Value = '  +9.12%  '
$ws.Range('E16').Value = '  +6.47%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '26.64'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +11.17%  '
$ws.Range('D18').Value = '2.689.60'
$ws.Range('E18').Value = '  +11.29%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '9.35'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +30.87%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.03'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +11.48%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '374.71'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +9.30%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.29'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +16.71%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.07'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +4.90%  '
$ws.Range('E24').Value = '  +4.44%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '70.32'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +6.18%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '4.18'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +9.88%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.42'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +10.81%  '
$ws.Range('D29').Value = '2.817.88'
$ws.Range('E29').Value = '  +10.22%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.00'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.36%  '
$ws.Range('D31').Value = '0.0₃0949'
$ws.Range('E31').Value = '  +11.79%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '522.57'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +14.95%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.41'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +13.40%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '7.76'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.96%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.76'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +9.44%  '
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '162.79'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.84%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.119'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +8.49%  '
$ws.Range('E39').Value = '  +5.96%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '19.40'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.54%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.02'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +14.21%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '170.77'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +26.41%  '
$ws.Range('E44').Value = '  +12.05%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.332'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +9.49%  '
$ws.Range('E46').Value = '  +9.96%  '
$ws.Range('E47').Value = '  +14.15%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '39.46'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +4.40%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0846'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +16.72%  '
$ws.Range('E50').Value = '  +7.83%  '
$ws.Range('E51').Value = '  +10.43%  '
